$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before current row 321 (Excel row index), shifting rows 321..432 down to 322..433
$ws.Rows.Item(321).Insert()

# Populate the newly inserted row 321 with the new data point
$ws.Cells.Item(321, 1).Value = 5
$ws.Cells.Item(321, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(321, 3).Value = "Maule"
$ws.Cells.Item(321, 4).Value = (Get-Date -Year 2021 -Month 11 -Day 24 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(321, 5).Value = 7
$ws.Cells.Item(321, 6).Value = "Fruta"
$ws.Cells.Item(321, 7).Value = 100108
$ws.Cells.Item(321, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(321, 9).Value = 100108006
$ws.Cells.Item(321, 10).Value = "Plátano"
$ws.Cells.Item(321, 11).Value = "Sin especificar"
$ws.Cells.Item(321, 12).Value = "Pintón"
$ws.Cells.Item(321, 13).Value = 520
$ws.Cells.Item(321, 14).Value = 23000
$ws.Cells.Item(321, 15).Value = 23000
$ws.Cells.Item(321, 16).Value = 23000
$ws.Cells.Item(321, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(321, 18).Value = "Ecuador"
$ws.Cells.Item(321, 19).Value = 1150
$ws.Cells.Item(321, 20).Value = 20
